$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the OTP value shown in B2 (shared string "ABC12" -> "ABC123")
$ws.Range("B2").Value = "ABC123"

# Move / record the active selection used by the last test run (B11 -> B6)
$ws.Range("B6").Select()
